# Update Cd86-Cd28 LR-pair table with new TPM-derived values.
# Adds an "ECs" sending-cluster block (4 new rows) and rewrites every
# numeric column (E:T) for the Inflammatory-Mac/Neutrophils/Resolving-Mac
# blocks too, per the refreshed NATMI run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Cd86", "Cd28", "Inflammatory-Mac", 1, 0.3333333333333333, 0.483405, 1.450215, 0.003094960828876145, 0.003094960828876144, 3, 1, 1.225147333333333, 3.675442, 0.2944933560673559, 0.2944933560673559, 0.5922423466700001, 5.330181120030001, 0.0009114454013927415, 0.0009114454013927413),
  @("ECs", "Cd86", "Cd28", "MuSCs", 1, 0.3333333333333333, 0.483405, 1.450215, 0.003094960828876145, 0.003094960828876144, 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.001282073881300198, 0.001282073881300198, 0.002578321135, 0.023204890215, 0.000003967968442349317, 0.000003967968442349316),
  @("ECs", "Cd86", "Cd28", "Neutrophils", 1, 0.3333333333333333, 0.483405, 1.450215, 0.003094960828876145, 0.003094960828876144, 3, 1, 0.4766303333333334, 1.429891, 0.1145694584217375, 0.1145694584217375, 0.2304054862850001, 2.073649376565001, 0.0003545879860008317, 0.0003545879860008316),
  @("ECs", "Cd86", "Cd28", "Resolving-Mac", 1, 0.3333333333333333, 0.483405, 1.450215, 0.003094960828876145, 0.003094960828876144, 3, 1, 2.453075333333333, 7.359226, 0.5896551116296064, 0.5896551116296064, 1.18582888151, 10.67245993359, 0.001824959473040222, 0.001824959473040222),
  @("Inflammatory-Mac", "Cd86", "Cd28", "Inflammatory-Mac", 3, 1, 105.7018663333333, 317.105599, 0.676747521934545, 0.6767475219345449, 3, 1, 1.225147333333333, 3.675442, 0.2944933560673559, 0.2944933560673559, 129.5003596666398, 1165.503236999758, 0.1992976489447707, 0.1992976489447707),
  @("Inflammatory-Mac", "Cd86", "Cd28", "MuSCs", 3, 1, 105.7018663333333, 317.105599, 0.676747521934545, 0.6767475219345449, 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.001282073881300198, 0.001282073881300198, 0.5637785210665556, 5.074006689599, 0.000867640322106913, 0.0008676403221069129),
  @("Inflammatory-Mac", "Cd86", "Cd28", "Neutrophils", 3, 1, 105.7018663333333, 317.105599, 0.676747521934545, 0.6767475219345449, 3, 1, 0.4766303333333334, 1.429891, 0.1145694584217375, 0.1145694584217375, 50.38071578441212, 453.4264420597091, 0.07753459707629373, 0.07753459707629373),
  @("Inflammatory-Mac", "Cd86", "Cd28", "Resolving-Mac", 3, 1, 105.7018663333333, 317.105599, 0.676747521934545, 0.6767475219345449, 3, 1, 2.453075333333333, 7.359226, 0.5896551116296064, 0.5896551116296064, 259.2946409895971, 2333.651768906374, 0.3990476355913736, 0.3990476355913736),
  @("Neutrophils", "Cd86", "Cd28", "Inflammatory-Mac", 3, 1, 18.83134933333333, 56.494048, 0.1205661682058513, 0.1205661682058513, 3, 1, 1.225147333333333, 3.675442, 0.2944933560673559, 0.2944933560673559, 23.07117741880178, 207.640596769216, 0.0355059355031225, 0.0355059355031225),
  @("Neutrophils", "Cd86", "Cd28", "MuSCs", 3, 1, 18.83134933333333, 56.494048, 0.1205661682058513, 0.1205661682058513, 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.001282073881300198, 0.001282073881300198, 0.1004401402275556, 0.9039612620480001, 0.0001545747352251683, 0.0001545747352251683),
  @("Neutrophils", "Cd86", "Cd28", "Neutrophils", 3, 1, 18.83134933333333, 56.494048, 0.1205661682058513, 0.1205661682058513, 3, 1, 0.4766303333333334, 1.429891, 0.1145694584217375, 0.1145694584217375, 8.975592309863112, 80.78033078876801, 0.01381320059532849, 0.01381320059532849),
  @("Neutrophils", "Cd86", "Cd28", "Resolving-Mac", 3, 1, 18.83134933333333, 56.494048, 0.1205661682058513, 0.1205661682058513, 3, 1, 2.453075333333333, 7.359226, 0.5896551116296064, 0.5896551116296064, 46.19471854298311, 415.752466886848, 0.07109245737217516, 0.07109245737217515),
  @("Resolving-Mac", "Cd86", "Cd28", "Inflammatory-Mac", 3, 1, 31.17437066666666, 93.523112, 0.1995913490307275, 0.1995913490307275, 3, 1, 1.225147333333333, 3.675442, 0.2944933560673559, 0.2944933560673559, 38.19319709061156, 343.738773815504, 0.05877832621806994, 0.05877832621806994),
  @("Resolving-Mac", "Cd86", "Cd28", "MuSCs", 3, 1, 31.17437066666666, 93.523112, 0.1995913490307275, 0.1995913490307275, 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.001282073881300198, 0.001282073881300198, 0.1662737016791111, 1.496463315112, 0.0002558908555257673, 0.0002558908555257673),
  @("Resolving-Mac", "Cd86", "Cd28", "Neutrophils", 3, 1, 31.17437066666666, 93.523112, 0.1995913490307275, 0.1995913490307275, 3, 1, 0.4766303333333334, 1.429891, 0.1145694584217375, 0.1145694584217375, 14.85865068231023, 133.727856140792, 0.02286707276411443, 0.02286707276411443),
  @("Resolving-Mac", "Cd86", "Cd28", "Resolving-Mac", 3, 1, 31.17437066666666, 93.523112, 0.1995913490307275, 0.1995913490307275, 3, 1, 2.453075333333333, 7.359226, 0.5896551116296064, 0.5896551116296064, 76.47307971459023, 688.257717431312, 0.1176900591930173, 0.1176900591930173)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  for ($j = 0; $j -lt 20; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $data[$i][$j]
  }
}